# Applies the commit: renames header columns on "Weekly Quantity" and
# "Monthly Trend", and adds a new "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- Rename the "Requested quantity" header on both existing sheets ---
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the other sheets (0.75"/0.75"/1"/1"/0.5"/0.5").
$wsForecast.PageSetup.LeftMargin   = 54
$wsForecast.PageSetup.RightMargin  = 54
$wsForecast.PageSetup.TopMargin    = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Pick up matching cell formatting (bold/centered/bordered header, date
# number format for column A) by copying it from the "Weekly Quantity"
# sheet before overwriting with the forecast's own values.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:D1"))
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A55"))

# --- Header row ---
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# --- Forecast data rows 2..55 ---
$data = @(
    @(44976.99999999999, 0, -143.685820292379, 121.1626061997484),
    @(44983.99999999999, 0, -137.9306810165527, 130.5241127184159),
    @(44997.99999999999, 0, -136.3489969243094, 133.3306754166084),
    @(45011.99999999999, 0, -136.8218028907774, 135.6977067496454),
    @(45018.99999999999, 2, -128.917116201413, 126.9949201635883),
    @(45039.99999999999, 8, -131.3881636974626, 136.443127352557),
    @(45060.99999999999, 13, -121.0279046608314, 157.1167921911986),
    @(45095.99999999999, 23, -106.077179365514, 161.8175301085614),
    @(45102.99999999999, 25, -109.5609296042904, 161.5315769915752),
    @(45109.99999999999, 26, -113.820410159511, 156.1301690654763),
    @(45116.99999999999, 28, -98.79582410551234, 162.1062760597801),
    @(45123.99999999999, 30, -104.9496178479217, 163.3519774908711),
    @(45130.99999999999, 32, -101.7434362267775, 173.5535083211913),
    @(45151.99999999999, 38, -101.4850385788039, 166.9881703798819),
    @(45165.99999999999, 42, -85.89475513126565, 172.0574031308441),
    @(45172.99999999999, 44, -91.57492263539848, 175.7849583992979),
    @(45179.99999999999, 45, -87.52453373741538, 180.0522339875342),
    @(45186.99999999999, 47, -91.58991660830836, 185.3983346923988),
    @(45193.99999999999, 49, -80.40862015471629, 182.1141546436436),
    @(45200.99999999999, 51, -87.91822988138627, 190.0231849344646),
    @(45207.99999999999, 53, -81.78078797296868, 184.8073213347397),
    @(45214.99999999999, 55, -74.88615376495331, 194.7727637248639),
    @(45249.99999999999, 64, -69.59977243359963, 202.7584907389961),
    @(45256.99999999999, 66, -76.53681933476746, 202.6720561163611),
    @(45270.99999999999, 70, -67.00529033496838, 200.2556380002407),
    @(45277.99999999999, 72, -54.52487699742729, 201.385813305074),
    @(45298.99999999999, 78, -52.29247618354722, 215.9429448965279),
    @(45305.99999999999, 79, -56.28659533886998, 221.5602375558782),
    @(45361.99999999999, 95, -39.80401371033671, 236.9285510343268),
    @(45375.99999999999, 98, -41.09485705543347, 234.8824355415596),
    @(45389.99999999999, 102, -31.99302453405263, 229.8099407927404),
    @(45396.99999999999, 104, -27.10161296212845, 246.489202986933),
    @(45403.99999999999, 106, -20.73624087408005, 242.7621144978532),
    @(45410.99999999999, 108, -17.7232786237657, 246.6225742279252),
    @(45417.99999999999, 110, -22.22966800920741, 245.9437121428142),
    @(45424.99999999999, 112, -18.64662665989543, 238.0860878488791),
    @(45431.99999999999, 114, -11.23642282271392, 246.9658978350363),
    @(45445.99999999999, 117, -15.02153811487949, 245.5856367129069),
    @(45452.99999999999, 119, -12.21787292261642, 249.9292583843132),
    @(45459.99999999999, 121, -7.104737039449127, 255.1174376372996),
    @(45466.99999999999, 123, -6.087256414101694, 247.3977915280856),
    @(45529.99999999999, 140, 2.850233442514471, 269.2463076585706),
    @(45550.99999999999, 146, 1.916966270336593, 274.6301359996179),
    @(45557.99999999999, 148, 8.025701135270504, 281.3628061409869),
    @(45578.99999999999, 153, 20.66384064172419, 287.7027005180479),
    @(45634.99999999999, 168, 38.35765317951277, 305.8939042052416),
    @(45641.99999999999, 170, 31.93346196076326, 304.8733155627628),
    @(45648.99999999999, 172, 39.16848477728107, 311.9277994944585),
    @(45655.99999999999, 174, 33.05881834740382, 311.7702955671117),
    @(45662.99999999999, 176, 54.95784802060108, 312.5571812783019),
    @(45669.99999999999, 178, 41.32583131534481, 310.8537253588199),
    @(45676.99999999999, 180, 51.79315072336006, 306.6860851068677),
    @(45683.99999999999, 182, 38.55281711158099, 316.5492513797394),
    @(45690.99999999999, 184, 53.78608290069745, 328.5522886837455)
)

$rowCount = $data.Count
$arr = New-Object 'object[,]' $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $data[$i]
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i,$j] = $r[$j]
    }
}

$destRange = $wsForecast.Range("A2").Resize($rowCount, 4)
$destRange.Value = $arr
